$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- 1. Update the "last updated" timestamp in the header cell ---
$ws.Range("A1").Value2 = "Datos actualizados a 10 de Abril de 2020 a las 10:52"

# --- 2. Update per-country daily statistics (Casos totales, Nuevos casos,
#         Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
#         Rows are addressed at their CURRENT (pre-sort) positions. ---

# Austria (row 19): Casos activos, Recuperados, Casos criticos, Muertes change
$ws.Cells.Item(19, 4).Value2 = 6064
$ws.Cells.Item(19, 5).Value2 = 6954
$ws.Cells.Item(19, 7).Value2 = 24
$ws.Cells.Item(19, 8).Value2 = 319

# Australia (row 27): Casos totales, Nuevos casos, Recuperados change
$ws.Cells.Item(27, 2).Value2 = 6203
$ws.Cells.Item(27, 3).Value2 = 51
$ws.Cells.Item(27, 5).Value2 = 3163

# Indonesia (row 40): full refresh, overtakes Mexico
$ws.Cells.Item(40, 2).Value2 = 3512
$ws.Cells.Item(40, 3).Value2 = 219
$ws.Cells.Item(40, 4).Value2 = 282
$ws.Cells.Item(40, 5).Value2 = 2924
$ws.Cells.Item(40, 6).Value2 = 0
$ws.Cells.Item(40, 7).Value2 = 26
$ws.Cells.Item(40, 8).Value2 = 306

# Estonia (row 65): full refresh, overtakes Irak
$ws.Cells.Item(65, 2).Value2 = 1258
$ws.Cells.Item(65, 3).Value2 = 51
$ws.Cells.Item(65, 4).Value2 = 93
$ws.Cells.Item(65, 5).Value2 = 1141
$ws.Cells.Item(65, 6).Value2 = 9
$ws.Cells.Item(65, 8).Value2 = 24

# Libano (row 85): Casos totales, Nuevos casos, Recuperados change
$ws.Cells.Item(85, 2).Value2 = 583
$ws.Cells.Item(85, 3).Value2 = 1
$ws.Cells.Item(85, 5).Value2 = 497

# Banglades (row 103): full refresh, overtakes Niger..San Marino
$ws.Cells.Item(103, 2).Value2 = 424
$ws.Cells.Item(103, 3).Value2 = 94
$ws.Cells.Item(103, 4).Value2 = 33
$ws.Cells.Item(103, 5).Value2 = 364
$ws.Cells.Item(103, 6).Value2 = 1
$ws.Cells.Item(103, 7).Value2 = 6
$ws.Cells.Item(103, 8).Value2 = 27

# Estado de Palestina (row 108): Casos totales, Nuevos casos, Recuperados change
$ws.Cells.Item(108, 2).Value2 = 266
$ws.Cells.Item(108, 3).Value2 = 3
$ws.Cells.Item(108, 5).Value2 = 221

# Sri Lanka (row 115): Casos activos, Recuperados change
$ws.Cells.Item(115, 4).Value2 = 50
$ws.Cells.Item(115, 5).Value2 = 133

# --- 3. Re-sort the country table by "Casos totales" (column B) descending,
#         same as the site does after every data refresh. ---
$sortRange = $ws.Range("A4:H216")
$sortKey = $ws.Range("B4:B216")
$sortRange.Sort($sortKey, 2)
